$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2

# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.3
$ws.Range("J6").Value = 2.5
$ws.Range("K6").Value = 1.95
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("U6").Value = 2.38
$ws.Range("V6").Value = 1.53
$ws.Range("AG6").Value = 10
$ws.Range("AL6").Value = 51
$ws.Range("AP6").Value = 26
$ws.Range("AQ6").Value = 34

# Row 10
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 2.88
$ws.Range("L10").Value = 3.4
$ws.Range("AK10").Value = 21
$ws.Range("AW10").Value = 5
$ws.Range("AX10").Value = 15
$ws.Range("BA10").Value = 51
